$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Re-point the three tables (slides 14, 15, 16) from the custom
#    "Table_0" style onto the built-in "No Style, Table Grid" style.
# ---------------------------------------------------------------------------
$newTableStyleId = "{BF65B155-BF0F-4922-90DA-8D23172560BA}"
foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyleId)
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Swap the theme color palette back to the stock Office theme colors
#    (the deck's master theme was carrying the "Integral / Red Violet"
#    palette; restore the default Office palette values).
# ---------------------------------------------------------------------------
$officeColors = @{
    1  = 0          # dk1      000000
    2  = 16777215   # lt1      FFFFFF
    3  = 6968388     # dk2      44546A
    4  = 15132391    # lt2      E7E6E6
    5  = 13998939    # accent1  5B9BD5
    6  = 3243501      # accent2  ED7D31
    7  = 10855845     # accent3  A5A5A5
    8  = 49407        # accent4  FFC000
    9  = 12874308     # accent5  4472C4
    10 = 4697456       # accent6  70AD47
    11 = 12673797      # hlink    0563C1
    12 = 7491477       # folHlink 954F72
}

$slideForTheme = $p.Slides.Item(1)
$themeColors = $slideForTheme.ThemeColorScheme
for ($idx = 1; $idx -le 12; $idx++) {
    $themeColors.Item($idx).RGB = $officeColors[$idx]
}
